# coverity_scan1.xlsx - "Add files via upload"
#
# Replaces the single data row's Coverity finding (function / issue /
# filename) and its line number, widens column A slightly, and moves the
# sheet's active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the finding reported in row 2 ---------------------------------
$ws.Range("A2").Value = "setDeviceInitializationContext"
$ws.Range("B2").Value = "COPY_INSTEAD_OF_MOVE"
$ws.Range("C2").Value = "entservices-softwareupdate/MaintenanceManager/MaintenanceManager.cpp "
$ws.Range("D2").Value = 1369

# --- Widen column A to fit the new (longer) function name -----------------
$ws.Columns.Item(1).ColumnWidth = 17.5

# --- Move the active selection to C9, as last left by the editor ----------
$ws.Range("C9").Select() | Out-Null
